# Update cryptocurrency price/volume snapshot (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.959.50"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.794.44"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.50"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5327"
$ws.Range("E7").Value = "  -2.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3967"
$ws.Range("E8").Value = "  +4.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07463"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.35"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.084"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.194"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.491"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.37"
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.792.64"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.41"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001060"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06573"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.955"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.983.16"
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.08"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.089"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.70"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.20"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.995.83"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.303"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.06"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1089"
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("E32").Value = "  -2.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.673"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.506"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07081"
$ws.Range("E35").Value = "  +5.29%  "
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.144"
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02274"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.397"
$ws.Range("E39").Value = "  -2.80%  "
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.188"
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6122"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.415"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.35"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.672"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5720"
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.09"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.919"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06805"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.27"
$ws.Range("E51").Value = "  -1.49%  "
